# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# on the active worksheet to the latest scraped values, matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Leading apostrophes are used on column D ("Price") assignments whose new
# text would otherwise be auto-parsed by Excel as a number (e.g. "1.001",
# "303.24"), so the cell keeps storing the literal text string exactly as
# scraped (mirrors the original inline-string cell content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.356.29'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.633.21'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '''303.24'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = '''0.3822'
$ws.Range("E7").Value = '  +0.96%  '
$ws.Range("D8").Value = '''51.99'
$ws.Range("D9").Value = '''0.3554'
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").Value = '''0.08130'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '''1.223'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = '''7.302'
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("D17").Value = '1.630.01'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '''94.65'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").Value = '''0.06940'
$ws.Range("D20").Value = '''6.576'
$ws.Range("E20").Value = '  +2.56%  '
$ws.Range("D21").Value = '''17.31'
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '''12.41'
$ws.Range("D24").Value = '23.349.12'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '''2.557'
$ws.Range("E25").Value = '  +4.52%  '
$ws.Range("D26").Value = '''3.121'
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").Value = '''20.99'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = '''151.41'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = '''5.262'
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("D30").Value = '''132.92'
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").Value = '1.809.46'
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").Value = '''2.147'
$ws.Range("E32").Value = '  -6.43%  '
$ws.Range("D33").Value = '''1.071'
$ws.Range("E33").Value = '  +13.25%  '
$ws.Range("D34").Value = '''6.500'
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("D35").Value = '''11.51'
$ws.Range("E35").Value = '  +5.30%  '
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("D37").Value = '''0.2485'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").Value = '''0.08735'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '''5.916'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").Value = '''0.06948'
$ws.Range("E40").Value = '  -2.23%  '
$ws.Range("D41").Value = '''0.6948'
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("D42").Value = '''1.323'
$ws.Range("E42").Value = '  -2.53%  '
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").Value = '''15.42'
$ws.Range("E44").Value = '  -4.27%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '''0.6380'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").Value = '''2.268'
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("D48").Value = '''3.950'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").Value = '''0.07920'
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = '''129.46'
$ws.Range("E50").Value = '  +3.32%  '
$ws.Range("D51").Value = '''1.182'
$ws.Range("E51").Value = '  -1.13%  '
